$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 2025-11-28 daily totals for both stations (rows 56-57).

# Row 56: 四方坪站
$ws.Cells.Item(56, 1).Value = 45989
$ws.Cells.Item(56, 2).Value = "四方坪站"
$ws.Cells.Item(56, 3).Value = 9150.9
$ws.Cells.Item(56, 4).Value = 8225.89
$ws.Cells.Item(56, 5).Value = 3056.78
$ws.Cells.Item(56, 6).Value = 401

# Row 57: 高岭站
$ws.Cells.Item(57, 1).Value = 45989
$ws.Cells.Item(57, 2).Value = "高岭站"
$ws.Cells.Item(57, 3).Value = 5278.48
$ws.Cells.Item(57, 4).Value = 4520.2
$ws.Cells.Item(57, 5).Value = 1352.8
$ws.Cells.Item(57, 6).Value = 185

$ws.Range("H56").Select()
